# Auto refresh - 16-02-2026 10:58:03.64
# Advances the "as of" dates one day forward (Today_Date / Last_Data_Till)
# on the Excel_vs_ML sheet, recomputes the pacing metrics that depend on
# Last_Data_Till for campaigns that are currently LIVE as of the new date,
# and stamps the new refresh timestamp on the Exec_Summary sheet.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Excel_vs_ML")
$sum = $wb.Worksheets.Item("Exec_Summary")

$newToday        = 46069   # 2026-02-16
$newLastDataTill = 46068   # 2026-02-15

# Recomputed Days_Elapsed / Days_Left / Expected_Spend_Till_Date /
# Pacing_%_vs_Ideal / Ideal_Daily_Spend for the rows whose flight window
# (Flight_Start_Date .. Flight_End_Date) contains the new Last_Data_Till.
$overrides = @{
    3  = @{ I=53;  J=3;  K=522991.35; L=16.79; N=154923.85 }
    8  = @{ I=25;  J=34; K=57983.44;  L=71;    N=2813.86 }
    10 = @{ I=22;  J=35; K=173243.5;  L=20.43; N=11813.39 }
    12 = @{ I=45;  J=37; K=307211.01; L=20.86; N=13397.97 }
    13 = @{ I=31;  J=50; K=64823.5;   L=59.72; N=2613.27 }
    20 = @{ I=46;  J=44; K=216547.51; L=44.46; N=7440.88 }
    27 = @{ I=59;  J=24; K=410731.61; L=26.11; N=19606.63 }
    32 = @{ I=71;  J=49; K=340054.66; L=32.23; N=9492.790000000001 }
    33 = @{ I=41;  J=25; K=80358.12;  L=55.77; N=3381.56 }
    37 = @{ I=3;   J=27; K=53402.45;  L=3.53;  N=19708.91 }
    40 = @{ I=46;  J=19; K=84986.38;  L=47.23; N=4207.71 }
    42 = @{ I=4;   J=43; K=48873.82;  L=5.52;  N=13292.29 }
    45 = @{ I=65;  J=14; K=291768.55; L=31.24; N=18819.13 }
}

for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 5).Value = $newToday         # E: Today_Date
    $ws.Cells.Item($row, 6).Value = $newLastDataTill   # F: Last_Data_Till

    if ($overrides.ContainsKey($row)) {
        $vals = $overrides[$row]
        $ws.Cells.Item($row, 9).Value  = $vals.I   # I: Days_Elapsed
        $ws.Cells.Item($row, 10).Value = $vals.J   # J: Days_Left
        $ws.Cells.Item($row, 11).Value = $vals.K   # K: Expected_Spend_Till_Date
        $ws.Cells.Item($row, 12).Value = $vals.L   # L: Pacing_%_vs_Ideal
        $ws.Cells.Item($row, 14).Value = $vals.N   # N: Ideal_Daily_Spend
    }
}

$sum.Range("B5").Value = "2026-02-16 05:28 UTC"
